$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (Late/Heading/Outstanding columns
# shift one to the right, from N/O/P to O/P/Q) to make room for the new
# "Variable Instalments" related column.
$ws.Columns("N:N").Insert() | Out-Null

# "Repayment Schedule" becomes the active sheet/tab, with cell S8 selected.
$ws.Activate()
$ws.Range("S8").Select() | Out-Null
